$wb = $excel.ActiveWorkbook

# --- Add the two new columns (Laps, WinnerId) to the RaceWeekends sheet,
#     ahead of the existing ConstructorId column ---
$race = $wb.Worksheets.Item("RaceWeekends")
$race.Range("D1").Value = "Laps"
$race.Range("E1").Value = "WinnerId"
$race.Range("F1").Value = "ConstructorId"

# --- Reorder the sheet tabs: move "Countries" so it sits right after "Tracks" ---
$countries = $wb.Worksheets.Item("Countries")
$tracks = $wb.Worksheets.Item("Tracks")
$countries.Move($null, $tracks)

# --- Update the RaceWeekends selection to A2:F38 ---
$race.Activate()
$race.Range("A2:F38").Select()

# --- Finally, make "Constructors" (the last tab) the active/selected sheet ---
$constructors = $wb.Worksheets.Item("Constructors")
$constructors.Activate()
